$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 data - new question/answer entry
# (order matters for shared-string table index assignment)
$ws.Range("B3").Value = "Calling the ELEV is done by Up/Down switches "
$ws.Range("D3").Value = "How many floors are there?"
$ws.Range("E3").Value = "Floors number doesn't matter, because the motor will simulate the direction without the need for any counting"
$ws.Range("F3").Value = "22/1/2020"
$ws.Range("G3").Value = "24/1/2020"
$ws.Range("C3").Value = "-"

# Column widths
$ws.Range("B1").ColumnWidth = 47.7109375
$ws.Range("D1").ColumnWidth = 27.140625
$ws.Range("E1").ColumnWidth = 106.5703125

# Sheet view: scroll + selection
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C4").Select()
